$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of time-record data (row 8), mirroring the existing rows.
$ws.Range("A8").Value = "25.11.2019"
$ws.Range("B8").Value = 0.60069444444444442
$ws.Range("C8").Value = 0.61111111111111105
$ws.Range("D8").Formula = "=C8-B8"
$ws.Range("D8").NumberFormat = "[$]hh:mm;@"
$ws.Range("E8").Value = "Calculator Control Unit"
$ws.Range("F8").Value = "Architecture, Testbench"

# Match formatting of the prior same-day rows (B6:C7 use style index 2 -> numFmtId 20).
$ws.Range("B8:C8").NumberFormat = "h:mm"

# Move the active selection down to the next empty row, like Excel would after
# entering data in the row above.
$ws.Range("A9").Select()
